# Apply the edits described by the diff:
#  - C11 changes from "Y" to "N"
#  - C13 changes from "N" to "Y"
#  - Active selection moves to E13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C11").Value = "N"
$ws.Range("C13").Value = "Y"

$ws.Range("E13").Select() | Out-Null
